# Generate Report for Handback
# Update handback status timestamps / status for the two most recently
# processed files (daa50092... row 4, fa04b78c... row 5) across the
# Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet - "Latest HO Xliff Generate Date" column (G)
$wsOverview.Range("G4").Value = "2016-08-24 06:16:45"
$wsOverview.Range("G5").Value = "2016-08-24 06:16:45"

# zh-cn sheet - Status (E), Correspond Handoff Datetime (H),
# Correspond Handback DateTime (K)
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-24 06:16:40"
$wsZhCn.Range("H5").Value = "2016-08-24 06:16:40"
$wsZhCn.Range("K4").Value = "2016-08-24 06:16:57"
$wsZhCn.Range("K5").Value = "2016-08-24 06:16:57"

# de-de sheet - Status (E), Correspond Handoff Datetime (H),
# Correspond Handback DateTime (K)
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-24 06:16:45"
$wsDeDe.Range("H5").Value = "2016-08-24 06:16:45"
$wsDeDe.Range("K4").Value = "2016-08-24 06:17:12"
$wsDeDe.Range("K5").Value = "2016-08-24 06:17:12"
